# Generate Report for Handback
#
# This applies the "handback" report-generation pass to the localization
# status workbook:
#   - Overview + per-language "Status" cells flip from "Ready for handoff"
#     to "Handed back: in sync with en-US"
#   - Each language sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for both data rows,
#     with hyperlinks added on the newly-populated "Latest Target File"
#     cells (mirroring the existing source-file hyperlinks in column A)
#   - A couple of columns get widened to fit the newly-written long file
#     names / timestamps

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a87fce5ff60baf04c9c25bb8b0c17ef4432fd25f/e2e/6d6bd283-653b-4d39-ab19-220919ac0371.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a87fce5ff60baf04c9c25bb8b0c17ef4432fd25f/e2e/6e36b6d6-9dbc-4318-80cd-90cda1757c77.md"
$mdName1 = "6d6bd283-653b-4d39-ab19-220919ac0371.md"
$mdName2 = "6e36b6d6-9dbc-4318-80cd-90cda1757c77.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

# Row 2 (6d6bd283...)
$wsZh.Range("I2").Value = $mdName1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, [System.Type]::Missing, [System.Type]::Missing, $mdName1)
$wsZh.Range("J2").Value = "6d6bd283-653b-4d39-ab19-220919ac0371.b18d37bb9665b241dd78eb399168221927b6eaf6.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 02:41:12"

# Row 3 (6e36b6d6...)
$wsZh.Range("I3").Value = $mdName2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, [System.Type]::Missing, [System.Type]::Missing, $mdName2)
$wsZh.Range("J3").Value = "6e36b6d6-9dbc-4318-80cd-90cda1757c77.62ccebb56ab3b2ce94297d7429b89c3ef5db44dc.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 02:41:12"

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.14
$wsZh.Columns.Item(10).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# Row 2 (6d6bd283...)
$wsDe.Range("I2").Value = $mdName1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, [System.Type]::Missing, [System.Type]::Missing, $mdName1)
$wsDe.Range("J2").Value = "6d6bd283-653b-4d39-ab19-220919ac0371.b18d37bb9665b241dd78eb399168221927b6eaf6.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 02:41:18"

# Row 3 (6e36b6d6...)
$wsDe.Range("I3").Value = $mdName2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, [System.Type]::Missing, [System.Type]::Missing, $mdName2)
$wsDe.Range("J3").Value = "6e36b6d6-9dbc-4318-80cd-90cda1757c77.62ccebb56ab3b2ce94297d7429b89c3ef5db44dc.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 02:41:18"

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.14
$wsDe.Columns.Item(10).ColumnWidth = 39.14
